{"js": "// Replace each three-digit-by-one-digit multiplication expression in the\n// document's table cells with its new value, per the commit's regenerated\n// problem set. Old -> new values are all unique, so a direct text search\n// and replace for each pair is safe (no collisions between old and new\n// sets).\nconst replacements = [\n  [\"355\u00d77=\", \"361\u00d73=\"],\n  [\"417\u00d76=\", \"869\u00d75=\"],\n  [\"888\u00d72=\", \"235\u00d79=\"],\n  [\"445\u00d73=\", \"316\u00d72=\"],\n  [\"899\u00d72=\", \"335\u00d74=\"],\n  [\"862\u00d74=\", \"560\u00d73=\"],\n  [\"540\u00d77=\", \"226\u00d72=\"],\n  [\"995\u00d75=\", \"982\u00d74=\"],\n  [\"902\u00d78=\", \"747\u00d75=\"],\n  [\"761\u00d79=\", \"601\u00d76=\"],\n  [\"560\u00d75=\", \"145\u00d78=\"],\n  [\"347\u00d74=\", \"430\u00d72=\"],\n  [\"939\u00d75=\", \"825\u00d74=\"],\n  [\"123\u00d76=\", \"537\u00d79=\"],\n  [\"322\u00d75=\", \"795\u00d79=\"],\n  [\"462\u00d79=\", \"179\u00d73=\"],\n  [\"245\u00d78=\", \"336\u00d74=\"],\n  [\"652\u00d74=\", \"769\u00d74=\"],\n  [\"839\u00d77=\", \"382\u00d77=\"],\n  [\"667\u00d76=\", \"954\u00d73=\"],\n  [\"103\u00d76=\", \"971\u00d72=\"],\n  [\"709\u00d74=\", \"862\u00d75=\"],\n  [\"274\u00d79=\", \"191\u00d72=\"],\n  [\"388\u00d72=\", \"759\u00d75=\"],\n  [\"933\u00d75=\", \"376\u00d78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication expression in the\n# document's table cells with its new value, per the commit's regenerated\n# problem set. Old -> new values are all unique, so a direct Find/Replace\n# for each pair is safe (no collisions between old and new sets).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"355\u00d77=\"; New = \"361\u00d73=\"},\n    @{Old = \"417\u00d76=\"; New = \"869\u00d75=\"},\n    @{Old = \"888\u00d72=\"; New = \"235\u00d79=\"},\n    @{Old = \"445\u00d73=\"; New = \"316\u00d72=\"},\n    @{Old = \"899\u00d72=\"; New = \"335\u00d74=\"},\n    @{Old = \"862\u00d74=\"; New = \"560\u00d73=\"},\n    @{Old = \"540\u00d77=\"; New = \"226\u00d72=\"},\n    @{Old = \"995\u00d75=\"; New = \"982\u00d74=\"},\n    @{Old = \"902\u00d78=\"; New = \"747\u00d75=\"},\n    @{Old = \"761\u00d79=\"; New = \"601\u00d76=\"},\n    @{Old = \"560\u00d75=\"; New = \"145\u00d78=\"},\n    @{Old = \"347\u00d74=\"; New = \"430\u00d72=\"},\n    @{Old = \"939\u00d75=\"; New = \"825\u00d74=\"},\n    @{Old = \"123\u00d76=\"; New = \"537\u00d79=\"},\n    @{Old = \"322\u00d75=\"; New = \"795\u00d79=\"},\n    @{Old = \"462\u00d79=\"; New = \"179\u00d73=\"},\n    @{Old = \"245\u00d78=\"; New = \"336\u00d74=\"},\n    @{Old = \"652\u00d74=\"; New = \"769\u00d74=\"},\n    @{Old = \"839\u00d77=\"; New = \"382\u00d77=\"},\n    @{Old = \"667\u00d76=\"; New = \"954\u00d73=\"},\n    @{Old = \"103\u00d76=\"; New = \"971\u00d72=\"},\n    @{Old = \"709\u00d74=\"; New = \"862\u00d75=\"},\n    @{Old = \"274\u00d79=\"; New = \"191\u00d72=\"},\n    @{Old = \"388\u00d72=\"; New = \"759\u00d75=\"},\n    @{Old = \"933\u00d75=\"; New = \"376\u00d78=\"}\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $pair.New, 2)\n}\n"}
